$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Add "Number of Past Year MedStar visits" row to the last
#    missingness-summary table, right before its "Total" row.
# ------------------------------------------------------------------
$t = $d.Tables.Item($d.Tables.Count)

$totalRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cellText = $t.Rows.Item($i).Cells.Item(1).Range.Text.TrimEnd([char]7).TrimEnd([char]13).Trim()
    if ($cellText -eq "Total") {
        $totalRow = $t.Rows.Item($i)
        break
    }
}

$newRow = $t.Rows.Add($totalRow)
$newRow.Cells.Item(1).Range.Text = "Number of Past Year MedStar visits"
$newRow.Cells.Item(2).Range.Text = "0"
$newRow.Cells.Item(3).Range.Text = "347"

# ------------------------------------------------------------------
# 2. Remove the trailing density-plot image paragraph (the one right
#    after the "Below are density plots ..." paragraph, at the very
#    end of the document body).
# ------------------------------------------------------------------
for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
    $shp = $d.InlineShapes.Item($i)
    if ($shp.Width -eq 360 -and $shp.Height -eq 288) {
        $para = $shp.Range.Paragraphs.Item(1)
        $para.Range.Delete()
    }
}
